# resultados atualizados. Closes #23.
#
# Applies the content updates from the "centro.spt" results table:
#   - n: 313 -> 312
#   - header label "Grupo (%)" -> "CAE (%)"
#   - p-value for that row: 0.583 -> 0.532
#   - removes the now-obsolete "CIRURGIA VASCULAR" data row
#   - JOELHO row count: 114 (36.4) -> 114 (36.5)
#   - QUADRIL row count: 71 (22.7) -> 71 (22.8)
#   - TRAUMA row count: 78 (24.9) -> 78 (25.0)

$d = $word.ActiveDocument

# --- simple text substitutions (each string is unique in the document) ---

$d.Content.Find.Execute("313", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "312", 2) | Out-Null

$d.Content.Find.Execute("Grupo (%)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CAE (%)", 2) | Out-Null

$d.Content.Find.Execute("0.583", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0.532", 2) | Out-Null

$d.Content.Find.Execute("114 (36.4)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "114 (36.5)", 2) | Out-Null

$d.Content.Find.Execute("71 (22.7)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "71 (22.8)", 2) | Out-Null

$d.Content.Find.Execute("78 (24.9)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "78 (25.0)", 2) | Out-Null

# --- remove the obsolete "CIRURGIA VASCULAR" row from the results table ---

$t = $d.Tables.Item(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "*CIRURGIA VASCULAR*") {
        $row.Delete()
        break
    }
}
